$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.408.90"
$ws.Range("E2").Value = "  +1.26%  "

$ws.Range("D3").Value = "2.379.90"
$ws.Range("E3").Value = "  +1.39%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'551.91"
$ws.Range("E5").Value = "  +1.68%  "

$ws.Range("E6").Value = "  +1.89%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "'0.524"
$ws.Range("E8").Value = "  +0.82%  "

$ws.Range("D9").Value = "2.380.37"
$ws.Range("E9").Value = "  +1.39%  "

$ws.Range("E10").Value = "  +4.39%  "

$ws.Range("E11").Value = "  +2.04%  "

$ws.Range("D12").Value = "'5.35"
$ws.Range("E12").Value = "  +2.70%  "

$ws.Range("E13").Value = "  +4.00%  "

$ws.Range("D14").Value = "'25.58"
$ws.Range("E14").Value = "  +3.76%  "

$ws.Range("E15").Value = "  +5.40%  "

$ws.Range("D16").Value = "2.808.72"
$ws.Range("E16").Value = "  +1.30%  "

$ws.Range("D17").Value = "61.264.04"
$ws.Range("E17").Value = "  +1.49%  "

$ws.Range("D18").Value = "2.381.48"
$ws.Range("E18").Value = "  +1.47%  "

$ws.Range("D19").Value = "'11.00"
$ws.Range("E19").Value = "  +4.19%  "

$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "'4.15"
$ws.Range("E20").Value = "  +2.48%  "

$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'321.46"
$ws.Range("E21").Value = "  +2.64%  "

$ws.Range("D22").Value = "'6.72"
$ws.Range("E22").Value = "  +3.01%  "

$ws.Range("E23").Value = "  +0.18%  "

$ws.Range("E24").Value = "  -4.78%  "

$ws.Range("E25").Value = "  +2.58%  "

$ws.Range("E26").Value = "  +8.92%  "

$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.25%  "

$ws.Range("D28").Value = "2.496.83"
$ws.Range("E28").Value = "  +1.34%  "

$ws.Range("D29").Value = "'8.21"
$ws.Range("E29").Value = "  +3.37%  "

$ws.Range("D30").Value = "'514.74"
$ws.Range("E30").Value = "  +2.79%  "

$ws.Range("D31").Value = "0.0₃0900"
$ws.Range("E31").Value = "  +1.65%  "

$ws.Range("E32").Value = "  +1.15%  "

$ws.Range("E33").Value = "  +4.10%  "

$ws.Range("D34").Value = "'1.84"
$ws.Range("E34").Value = "  +3.09%  "

$ws.Range("E35").Value = "  +1.22%  "

$ws.Range("E36").Value = "  +0.03%  "

$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "'4.70"
$ws.Range("E37").Value = "  +4.05%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").Value = "'5.51"
$ws.Range("E38").Value = "  +5.70%  "

$ws.Range("D39").Value = "'1.89"
$ws.Range("E39").Value = "  +6.72%  "

$ws.Range("E40").Value = "  +2.30%  "

$ws.Range("D41").Value = "'18.54"
$ws.Range("E41").Value = "  +1.35%  "

$ws.Range("D42").Value = "'146.68"
$ws.Range("E42").Value = "  +6.10%  "

$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("D44").Value = "'41.30"
$ws.Range("E44").Value = "  +3.18%  "

$ws.Range("D45").Value = "'148.11"
$ws.Range("E45").Value = "  +6.82%  "

$ws.Range("D46").Value = "'2.15"
$ws.Range("E46").Value = "  +3.26%  "

$ws.Range("E47").Value = "  +2.53%  "

$ws.Range("D48").Value = "'0.0528"
$ws.Range("E48").Value = "  +3.54%  "

$ws.Range("D49").Value = "'19.65"
$ws.Range("E49").Value = "  +1.09%  "

$ws.Range("D50").Value = "'0.581"
$ws.Range("E50").Value = "  +2.57%  "

$ws.Range("E51").Value = "  +1.53%  "
